$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2
$ws.Range("A2").Value = "Potato Russel"
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = "Restaurant Depot"

# Add new row 3
$ws.Range("A3").Value = "Frozen Peas and carrot cut"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "Restaurant Depot"

# Add new row 4
$ws.Range("A4").Value = "Paneer - Not Appel, not Nanak"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "Restaurant Depot"
